# Applies the edits described by the commit "Ajout des scripts pour les pièces"
# to the journal de travail workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# Row 52: add a hyperlink in I52 pointing to the Unity documentation page
# used while implementing prefab instantiation.
# ---------------------------------------------------------------------------
$i52 = $ws.Range("I52")
$i52.Value = "https://docs.unity3d.com/Manual/InstantiatingPrefabs.html"
$ws.Hyperlinks.Add($i52, "https://docs.unity3d.com/Manual/InstantiatingPrefabs.html") | Out-Null
# Re-copy the existing hyperlink cell style (used by I17/I18/I22/I26/I35) onto
# I52 so the new cell matches the sheet's existing hyperlink formatting
# instead of the ad-hoc style Hyperlinks.Add creates.
$ws.Range("I17").Copy($i52)
$i52.Value = "https://docs.unity3d.com/Manual/InstantiatingPrefabs.html"

# ---------------------------------------------------------------------------
# Row 54: new work session (15:20 - 16:55) on 2024-05-07 (serial 45419).
# Only the Date/Début/Fin columns are filled in; the Durée formula already
# exists in the sheet and recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("B52:D52").Copy($ws.Range("B54:D54"))
$ws.Range("B54").Value = 45419
$ws.Range("C54").Value = 0.63888888888888895
$ws.Range("D54").Value = 0.70486111111111116

# ---------------------------------------------------------------------------
# Row 55: only a date (2024-05-13, serial 45425) — this value used to sit on
# row 56 and is shifted up to row 55 to make room for the new row 54 entry.
# ---------------------------------------------------------------------------
$ws.Range("B52").Copy($ws.Range("B55"))
$ws.Range("B55").Value = 45425

# Row 56 loses the date value that moved to row 55 (its Durée formula/cell
# stays in place).
$ws.Range("B56").Clear()

# ---------------------------------------------------------------------------
# Row 58: work session (08:50 - 09:25) on 2024-05-16 (serial 45428),
# subject "Documentation", description "Rédaction de la planification".
# ---------------------------------------------------------------------------
$ws.Range("B52:D52").Copy($ws.Range("B58:D58"))
$ws.Range("B58").Value = 45428
$ws.Range("C58").Value = 0.36805555555555558
$ws.Range("D58").Value = 0.3923611111111111

$ws.Range("F53").Copy($ws.Range("F58"))
$ws.Range("F58").Value = "Documentation"

$ws.Range("G52").Copy($ws.Range("G58"))
$ws.Range("G58").Value = "Rédaction de la planification"

# ---------------------------------------------------------------------------
# Row 59: work session (09:25 - 09:35) on 2024-05-16 (serial 45428),
# subject "Implémentation", description "Création des scripts pour les
# pièces".
# ---------------------------------------------------------------------------
$ws.Range("B52:D52").Copy($ws.Range("B59:D59"))
$ws.Range("B59").Value = 45428
$ws.Range("C59").Value = 0.3923611111111111
$ws.Range("D59").Value = 0.39930555555555558

$ws.Range("F52").Copy($ws.Range("F59"))
$ws.Range("F59").Value = "Implémentation"

$ws.Range("G52").Copy($ws.Range("G59"))
$ws.Range("G59").Value = "Création des scripts pour les pièces"

# ---------------------------------------------------------------------------
# Recalculate everything (Durée formulas on Feuil1, SUMIF/percentage
# formulas on Sheet1, and the chart that is bound to them).
# ---------------------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------------
# Update the window/selection state to match where the author left off.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G59").Select()

$wb.Save()
